# Applies the cryptos-list price/volume refresh described by the commit diff.
# D-column cells whose new text looks like a plain decimal number (e.g. "585.61")
# are forced to Text format first, so Excel keeps them as strings instead of
# auto-converting them to numeric values (matching the original inlineStr cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.569.52"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "3.426.83"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.90%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "3.423.58"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  +2.32%  "
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.415"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.68%  "
$ws.Range("D13").Value = "4.020.36"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("D16").Value = "66.480.84"
$ws.Range("E16").Value = "  +0.65%  "
$ws.Range("E17").Value = "  +1.98%  "
$ws.Range("D18").Value = "3.420.51"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "369.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.15%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000127"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.535"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.03%  "
$ws.Range("E27").Value = "  +2.64%  "
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("E36").Value = "  +1.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.867"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.02%  "
$ws.Range("E40").Value = "  +3.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.17%  "
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("D43").Value = "2.718.49"
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0690"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "335.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.97%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.45%  "
$ws.Range("E51").Value = "  +3.81%  "
